$wb = $excel.ActiveWorkbook

# Update OFF sheet (Wild Card round values for Road "R" row)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 176
$wsOff.Range("C3").Value = 108
$wsOff.Range("D3").Value = 51
$wsOff.Range("E3").Value = 22

# Update DEF sheet (Wild Card round values for Road "R" row)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 167
$wsDef.Range("C3").Value = 106
$wsDef.Range("D3").Value = 51
$wsDef.Range("E3").Value = 20
$wsDef.Range("G3").Value = 4
